# Edit script: add "2022-Q1" fund-holdings sheet before the "总计" summary sheet,
# and prepend a 2022-Q1 row to the "总计" summary table.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Turn the existing "总计" sheet into the new "2022-Q1" sheet (this
#    preserves its sheetId/position slot) and add a brand-new "总计"
#    sheet right after it - this reproduces the sheetId renumbering
#    seen in the target workbook (old sheetId 5 -> new sheet, 总计 -> 6).
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Name = "2022-Q1"

$newTotalSheet = $wb.Worksheets.Add($null, $totalSheet)
$newTotalSheet.Name = "总计"
# Match the page margins used throughout the rest of the workbook
# (0.75in/0.75in/1in/1in/0.5in/0.5in, expressed here in points).
$newTotalSheet.PageSetup.LeftMargin = 54
$newTotalSheet.PageSetup.RightMargin = 54
$newTotalSheet.PageSetup.TopMargin = 72
$newTotalSheet.PageSetup.BottomMargin = 72
$newTotalSheet.PageSetup.HeaderMargin = 36
$newTotalSheet.PageSetup.FooterMargin = 36

# ------------------------------------------------------------------
# 2. Populate the "2022-Q1" sheet with the fund-holdings table.
#    "2021-Q4" has an identical column layout (B..H headers + index
#    column A), so it is used as the formatting template for the bold
#    / bordered / centered "style 2" header row and index column.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$template = $wb.Worksheets.Item("2021-Q4")
$q1.Cells.Clear()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $q1.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
}

$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q1.Range("A2:A35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Columns B, D, E, F, G hold numeric-looking strings that must stay text
# (fund codes with leading zeros, percentages kept at fixed precision, etc.)
$q1.Range("B2:B35").NumberFormat = "@"
$q1.Range("D2:G35").NumberFormat = "@"

$fundRows = @(
    @(0, "510900", "易方达恒生国企(QDII-ETF)", "101.39", "97.12", "4.56", "4.6234", 6),
    @(1, "159920", "华夏恒生ETF(QDII)", "151.31", "95.19", "2.75", "4.1610", 9),
    @(2, "501050", "华夏沪港通上证50AH优选指数（LOF）A", "25.94", "92.28", "7.23", "1.8755", 3),
    @(3, "007110", "国投瑞银港股通价值发现混合", "23.33", "93.33", "5.18", "1.2085", 6),
    @(4, "513550", "华泰柏瑞中证港股通50ETF", "31.28", "98.89", "3.83", "1.1980", 8),
    @(5, "009007", "兴全沪港深两年持有期混合", "24.80", "90.13", "3.14", "0.7787", 7),
    @(6, "513660", "华夏沪港通恒生ETF", "19.61", "97.34", "2.97", "0.5824", 9),
    @(7, "010010", "国投瑞银港股通6个月定期开放股票", "8.09", "93.58", "5.75", "0.4652", 6),
    @(8, "159850", "华夏恒生中国企业ETF（QDII）", "6.11", "93.95", "4.40", "0.2688", 6),
    @(9, "159960", "平安港股通恒生中国企业ETF", "4.08", "96.86", "5.16", "0.2105", 5),
    @(10, "501301", "华宝港股通恒生中国(香港上市)25指数(LOF)A", "2.82", "94.73", "6.73", "0.1898", 5),
    @(11, "513600", "南方恒生ETF", "5.89", "99.00", "3.02", "0.1779", 8),
    @(12, "159954", "南方恒生中国企业ETF", "3.35", "103.89", "4.85", "0.1625", 5),
    @(13, "160717", "嘉实恒生中国企业指数(QDII-LOF)", "2.57", "94.76", "4.44", "0.1141", 6),
    @(14, "161831", "银华恒生国企指数（QDII-LOF）", "2.29", "86.34", "4.79", "0.1097", 4),
    @(15, "007354", "创金合信港股通量化股票A", "3.84", "91.20", "2.13", "0.0818", 8),
    @(16, "164705", "汇添富恒生指数（QDII-LOF）A", "2.96", "92.23", "2.67", "0.0790", 9),
    @(17, "006355", "华宝港股通恒生中国(香港上市)25指数(LOF)C", "1.02", "94.73", "6.73", "0.0686", 5),
    @(18, "159712", "国泰中证港股通50ETF", "0.85", "95.21", "3.95", "0.0336", 7),
    @(19, "008407", "恒生前海恒生沪深港通细分行业龙头指数A", "0.37", "93.40", "7.14", "0.0264", 3),
    @(20, "160924", "大成恒生指数（QDII-LOF）", "0.89", "93.20", "2.73", "0.0243", 10),
    @(21, "513990", "招商上证港股通ETF", "0.59", "96.48", "3.69", "0.0218", 6),
    @(22, "006395", "华夏沪港通上证50AH优选指数（LOF）C", "0.25", "92.28", "7.23", "0.0181", 3),
    @(23, "513680", "建信港股通恒生中国企业ETF", "0.28", "96.77", "5.45", "0.0153", 6),
    @(24, "501067", "招商富时中国A-H50指数（LOF）A", "0.21", "94.63", "7.23", "0.0152", 2),
    @(25, "010789", "汇添富恒生指数（QDII-LOF）C", "0.37", "92.23", "2.67", "0.0099", 9),
    @(26, "159978", "建信中证沪港深粤港澳大湾区发展主题ETF", "0.31", "96.05", "3.05", "0.0095", 7),
    @(27, "001942", "前海开源沪港深汇鑫灵活配置混合A", "0.10", "90.39", "7.04", "0.0070", 8),
    @(28, "501309", "国泰恒生港股通指数（LOF）", "0.36", "92.35", "1.83", "0.0066", 9),
    @(29, "008408", "恒生前海恒生沪深港通细分行业龙头指数C", "0.08", "93.40", "7.14", "0.0057", 3),
    @(30, "001943", "前海开源沪港深汇鑫灵活配置混合C", "0.08", "90.39", "7.04", "0.0056", 8),
    @(31, "007357", "创金合信港股通量化股票C", "0.26", "91.20", "2.13", "0.0055", 8),
    @(32, "501068", "招商富时中国A-H50指数（LOF）C", "0.05", "94.63", "7.23", "0.0036", 2),
    @(33, "167302", "方正富邦恒生沪深港通大湾区综合指数（LOF）", "0.07", "93.79", "3.26", "0.0023", 4)
)

foreach ($row in $fundRows) {
    $r = [int]$row[0] + 2
    $q1.Cells.Item($r, 1).Value = [int]$row[0]
    $q1.Cells.Item($r, 2).Value = [string]$row[1]
    $q1.Cells.Item($r, 3).Value = [string]$row[2]
    $q1.Cells.Item($r, 4).Value = [string]$row[3]
    $q1.Cells.Item($r, 5).Value = [string]$row[4]
    $q1.Cells.Item($r, 6).Value = [string]$row[5]
    $q1.Cells.Item($r, 7).Value = [string]$row[6]
    $q1.Cells.Item($r, 8).Value = [int]$row[7]
}

# ------------------------------------------------------------------
# 3. Populate the new "总计" sheet with the summary table, the
#    2022-Q1 row now on top followed by the previously-existing rows.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$total.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summaryRows = @(
    @(0, "2022-Q1", 34, 16.57),
    @(1, "2021-Q4", 37, 13.37),
    @(2, "2021-Q2", 39, 20.39),
    @(3, "2021-Q1", 65, 45.66),
    @(4, "2020-Q4", 73, 40.03)
)

foreach ($row in $summaryRows) {
    $r = [int]$row[0] + 2
    $total.Cells.Item($r, 1).Value = [int]$row[0]
    $total.Cells.Item($r, 2).Value = [string]$row[1]
    $total.Cells.Item($r, 3).Value = [int]$row[2]
    $total.Cells.Item($r, 4).Value = [double]$row[3]
}

# Leave the window focused the way it started (no active-sheet change
# is part of the published diff).
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Host "done"
